$wb = $excel.ActiveWorkbook

# --- Update "Hoja1" conversion text (A1) ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$wsHoja1.Range("A1").Value = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 3.97 = 15284.77 pesos`n✅ 15284.77 pesos = 3.95 = 949.34 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

# --- Update "tasas" sheet rate cells ---
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("N10").Value = 251.623
$wsTasas.Range("O10").Value = 3846
$wsTasas.Range("N12").Value = 3869.9
$wsTasas.Range("O12").Value = 240.36
